# Apply the "gh-pages output regenerated" update to the 江西-漫展信息 workbook.
# Two worksheets contain the same underlying rows (展览 = exhibitions only,
# 全部类型 = all event types combined) so every numeric "想去人数" (want-to-go
# count) bump and the cancelled-event update must be mirrored on both sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" --------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F5").Value  = 4542
$ws1.Range("F7").Value  = 124
$ws1.Range("F9").Value  = 3065
$ws1.Range("F12").Value = 240
$ws1.Range("F13").Value = 588
$ws1.Range("F15").Value = 508
$ws1.Range("F16").Value = 355
$ws1.Range("F19").Value = 1297
$ws1.Range("F20").Value = 113
$ws1.Range("F21").Value = 1541
$ws1.Range("F22").Value = 124
$ws1.Range("F27").Value = 38

# Row 29: "上饶·次元重现夏日嘉年华" event got cancelled.
$ws1.Range("C29").Value = "上饶·次元重现夏日嘉年华（取消）"
$ws1.Range("G29").Value = "不可售"

$ws1.Range("F31").Value = 3477
$ws1.Range("G31").Value = 58.5
$ws1.Range("F32").Value = 740
$ws1.Range("F34").Value = 236
$ws1.Range("F36").Value = 1694

# ---- Sheet "全部类型" -----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value  = 4542
$ws4.Range("F7").Value  = 124
$ws4.Range("F9").Value  = 3065
$ws4.Range("F12").Value = 240
$ws4.Range("F13").Value = 588
$ws4.Range("F15").Value = 508
$ws4.Range("F17").Value = 355
$ws4.Range("F20").Value = 1297
$ws4.Range("F21").Value = 113
$ws4.Range("F22").Value = 1541
$ws4.Range("F23").Value = 124
$ws4.Range("F28").Value = 38

# Row 30: same cancelled event as above, mirrored here.
$ws4.Range("C30").Value = "上饶·次元重现夏日嘉年华（取消）"
$ws4.Range("G30").Value = "不可售"

$ws4.Range("F32").Value = 3477
$ws4.Range("G32").Value = 58.5
$ws4.Range("F34").Value = 740
$ws4.Range("F36").Value = 236
$ws4.Range("F38").Value = 1694
